$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns I and J, matching the formatting
# already used by the other header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Data rows for columns I (I0) and J (IF)
$data = @(
    @(2,5,5),
    @(3,6,6),
    @(4,8,8),
    @(5,5,5),
    @(6,8,8),
    @(7,7,7),
    @(8,10,11),
    @(9,7,7),
    @(10,6,7),
    @(11,6,7),
    @(12,9,9),
    @(13,7,8),
    @(14,8,8),
    @(15,6,6),
    @(16,8,8),
    @(17,9,9),
    @(18,9,9),
    @(19,5,6),
    @(20,8,8),
    @(21,7,7),
    @(22,8,8),
    @(23,8,9),
    @(24,6,7),
    @(25,7,7),
    @(26,8,8),
    @(27,6,6),
    @(28,8,8),
    @(29,8,8),
    @(30,6,7),
    @(31,9,9),
    @(32,8,8),
    @(33,8,8),
    @(34,6,6),
    @(35,8,9),
    @(36,3,4),
    @(37,8,8),
    @(38,3,3),
    @(39,9,9),
    @(40,6,6),
    @(41,6,6),
    @(42,5,5),
    @(43,8,8),
    @(44,7,7),
    @(45,8,8),
    @(46,7,7),
    @(47,8,8),
    @(48,5,6),
    @(49,7,7),
    @(50,7,7),
    @(51,9,9),
    @(52,6,7),
    @(53,5,6),
    @(54,4,4),
    @(55,7,8),
    @(56,7,7),
    @(57,8,8),
    @(58,7,7),
    @(59,9,9),
    @(60,5,5),
    @(61,7,8),
    @(62,5,5),
    @(63,5,6),
    @(64,6,6),
    @(65,8,8),
    @(66,7,7),
    @(67,4,4),
    @(68,3,3)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
